# Add the supervisor/cat-name note below the title, separated by two
# blank lines, without disturbing the existing "Git Hands On" title
# paragraph (Title style, centered).

$d = $word.ActiveDocument

# Append a paragraph break after the last paragraph in the document
# (the "Git Hands On" title). Re-resolve Content.End fresh each time so
# each new break is inserted after the one before it, rather than all
# being inserted at the same spot.
$d.Content.InsertParagraphAfter()
$d.Range($d.Content.End, $d.Content.End).InsertParagraphAfter()
$d.Range($d.Content.End, $d.Content.End).InsertParagraphAfter()

# The two new blank paragraphs (#2 and #3) inherited the Title/centered
# formatting from the paragraph mark they were split from. Word would
# normally fall back to the Title style's linked "next" style (Normal)
# for a freshly-typed paragraph; reproduce that explicitly so the blank
# paragraphs -- and, by extension, the paragraph that follows them --
# come out as plain Normal/left-aligned paragraphs, matching what a user
# typing after the title would get.
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $para.Range.ParagraphFormat.Style = "Normal"
    $para.Range.ParagraphFormat.Alignment = 0
}

# Now add the supervisor/cat text as a new 4th paragraph; it will pick
# up the already-reset Normal formatting from paragraph 3's mark.
$d.Range($d.Content.End, $d.Content.End).InsertAfter("Supervisor: This is Ousmane. Cats name Luna.")
